$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert two new rows before the existing data row (row 2), shifting the
# existing data down to row 4. Newly inserted rows inherit the header's
# format, so clear that before applying the real body formatting below.
$ws.Rows(2).Insert()
$ws.Rows(2).Insert()
$ws.Range("A2:E3").ClearFormats()
$lo.Resize($ws.Range("A1:E4"))

# New row 2 data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Extend the project scope to the team of 5"
$ws.Range("C2").Value = "Vikas, Sai Krishna, Revanth, Siri, Sai Teja"
$ws.Range("D2").Value = 43349
$ws.Range("E2").Value = "Closed"

# New row 3 data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Complete team deliverable 1 documents"
$ws.Range("C3").Value = "Vikas, Sai Krishna, Revanth, Siri, Sai Teja"
$ws.Range("D3").Value = 43355
$ws.Range("E3").Value = "Open"

# Existing row, now row 4: update S. No from 1 to 3
$ws.Range("A4").Value = 3

# Column C is now wide enough to show the full owner names
$ws.Columns(3).ColumnWidth = 36.17

# Body formatting: center alignment on all data rows, date format + center
# alignment on the Deadline column
$ws.Range("A2:E4").HorizontalAlignment = -4108
$ws.Range("D2:D4").NumberFormat = "m/d/yyyy"
